$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (Price / Volume(1h) columns).
# These columns are stored as text (e.g. "311.99", "-4.54%"), so each
# cell is forced to Text format before the write and the original
# (unstyled) cell style is restored afterward, to avoid Excel
# auto-converting the text into a number/percentage.
$updates = [ordered]@{
    "D2" = "311.99"
    "E2" = "-4.54%"
    "D3" = "40.64"
    "E3" = "-8.05%"
    "D4" = "5.118"
    "E4" = "-2.38%"
    "D5" = "0.07852"
    "E5" = "-5.69%"
    "D6" = "4.340"
    "E6" = "-2.12%"
    "D7" = "1.674"
    "E7" = "-13.53%"
    "D8" = "0.9242"
    "E8" = "-4.68%"
    "D9" = "0.1076"
    "E9" = "-4.41%"
    "D10" = "0.1789"
    "E10" = "-5.20%"
    "D11" = "0.09095"
    "E11" = "-5.13%"
    "D12" = "0.04450"
    "E12" = "-3.83%"
    "D13" = "7.195"
    "E13" = "-17.53%"
    "D14" = "0.1058"
    "E14" = "-0.01%"
    "D15" = "0.001265"
    "E15" = "-1.82%"
    "D16" = "0.005890"
    "E16" = "-3.23%"
    "E17" = "-1.31%"
    "D18" = "2.554"
    "E18" = "1.06%"
    "D19" = "0.3318"
    "E19" = "-0.85%"
    "D20" = "0.1381"
    "E20" = "0.62%"
    "D21" = "0.2655"
    "E21" = "2.82%"
    "E22" = "0.91%"
    "D23" = "0.001248"
    "E23" = "1.16%"
    "D24" = "0.004146"
    "E24" = "-6.56%"
    "D25" = "0.0001233"
    "E25" = "-5.33%"
    "D26" = "0.0003006"
    "E26" = "0.61%"
    "D38" = "0.02446"
    "E38" = "-8.71%"
    "D39" = "0.05284"
    "E39" = "-4.79%"
    "D40" = "0.008061"
    "E40" = "3.08%"
    "D41" = "0.1355"
    "E41" = "-3.71%"
    "D42" = "0.006563"
    "E42" = "-10.70%"
    "D43" = "0.002037"
    "E43" = "-4.22%"
    "D44" = "0.008275"
    "E44" = "5.26%"
    "D45" = "0.3109"
    "E45" = "-11.16%"
    "D46" = "0.00006805"
    "E46" = "-0.66%"
    "D47" = "0.00000000757"
    "E47" = "0.75%"
    "D48" = "0.003432"
    "E48" = "-1.67%"
    "D49" = "0.004138"
    "E49" = "16.85%"
    "D50" = "0.00002119"
    "E50" = "0.75%"
    "D51" = "0.0002018"
    "E51" = "0.75%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = $origStyle
}
